$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A (shifts existing Valid/T/Z/p-value data from A:D to B:E,
# preserving each column's original width).
$ws.Columns.Item(1).Insert()

# New column A holds the long "metric & metric" labels - widen it.
$ws.Columns.Item(1).ColumnWidth = 53.6

# Insert a new header row 1 (shifts the 13 data rows from 1:13 down to 2:14).
$ws.Rows.Item(1).Insert()

# Header row.
$ws.Range("B1").Value = "Valid"
$ws.Range("C1").Value = "T"
$ws.Range("D1").Value = "Z"
$ws.Range("E1").Value = "p-value"

# Row labels for the 13 data rows now sitting in rows 2-14.
$labels = @(
    "CyclomaticComplexity(CC) & CyclomaticComplexity(CC)",
    "MaintainabilityIndex & MaintainabilityIndex",
    "NbOperands & NbOperands",
    "NbOperands & EffortToImplement",
    "NbUniqueOperators & NbUniqueOperators",
    "NbOperators & NbOperators",
    "ProgramLength & ProgramLength",
    "VocabularySize & VocabularySize",
    "ProgramVolume & ProgramVolume",
    "DifficultyLevel & DifficultyLevel",
    "ProgramLevel & ProgramLevel",
    "EffortToImplement & NbOperands",
    "EffortToImplement & EffortToImplement"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $labels[$i]
}
